$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new twig template row (A5) with the extended twig syntax shared string.
# Setting Value2 on the new cell below the used range both creates the new
# shared string entry and extends the sheet dimension/used range to A1:H5,
# inheriting the column's default style (xf index 1 = Times New Roman 10pt)
# just like the other data cells in this column.
$ws.Range("A5").Value2 = '${twig: {% if record.bday %} {{record.bday|date("m/d/Y")}}{% endif %} }'

# Move/record the active selection onto the newly added cell.
$ws.Range("A5").Select()
